# Update existing cell values (A1:B3) and add two new rows (4 and 5),
# extending the data range from A1:B3 to A1:B5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.056297471750357576
$ws.Range("B1").Value = 0.056297471719789563

$ws.Range("A2").Value = 0.0091274790796969495
$ws.Range("B2").Value = -0.0091274791439819852

$ws.Range("A3").Value = -0.026666237263525513
$ws.Range("B3").Value = 0.026666237222573691

$ws.Range("A4").Value = -0.016135793659887369
$ws.Range("B4").Value = 0.016135793603759913

$ws.Range("A5").Value = 0.057438407269169753
$ws.Range("B5").Value = -0.057438407332522118

# Widen the two columns (A: 13.7109375 -> 14.42578125, B: 14.42578125 -> 15.42578125
# stored "characters" width). The host rounds ColumnWidth to 1/6-character pixel
# steps, so we pick the input that lands on the closest achievable stored width.
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
